$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header C1: "audioFalse" -> "currentPhase"
$ws.Range("C1").Value = "currentPhase"

# Column C (formerly per-row audio filenames) now holds a single shared
# value "train2P2" for both data rows.
$ws.Range("C2").Value = "train2P2"
$ws.Range("C3").Value = "train2P2"
